$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the hours value for row 12 (cascades running totals via formulas)
$ws.Range("B12").Value = 9

# Update the task description text in D12
$ws.Range("D12").Value = "Fixed vue routing; started to implement data receiving"

# Update the current selection to B13
$ws.Range("B13").Select()
